# Updated: po 21. 06. 2021
# Applies the Slovakia Covid daily-stats update:
#  - corrects several previously reported AgTests/AgPosit (columns F/G) values
#  - appends three new daily rows (471-473) for 2021-06-18 .. 2021-06-20

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to existing rows (columns F = AgTests, G = AgPosit) ---
$ws.Range("F308").Value = 16056

$ws.Range("F327").Value = 224078
$ws.Range("G327").Value = 2712

$ws.Range("F328").Value = 180331
$ws.Range("G328").Value = 2657

$ws.Range("F329").Value = 73163
$ws.Range("G329").Value = 1718

$ws.Range("F349").Value = 159670

$ws.Range("F359").Value = 320856
$ws.Range("G359").Value = 3333

$ws.Range("F362").Value = 228991

$ws.Range("F369").Value = 235137

$ws.Range("F377").Value = 176657
$ws.Range("G377").Value = 1824

$ws.Range("F379").Value = 180505

$ws.Range("F381").Value = 746942

$ws.Range("F387").Value = 351294
$ws.Range("G387").Value = 1663

$ws.Range("F388").Value = 730382

$ws.Range("F398").Value = 298822

$ws.Range("F407").Value = 158077
$ws.Range("G407").Value = 674

$ws.Range("F412").Value = 176418
$ws.Range("G412").Value = 647

$ws.Range("F413").Value = 149579

$ws.Range("F414").Value = 149059

$ws.Range("F419").Value = 149457

$ws.Range("F420").Value = 138754

$ws.Range("F421").Value = 153011

$ws.Range("F422").Value = 298118

$ws.Range("F431").Value = 166668
$ws.Range("G431").Value = 380

$ws.Range("F436").Value = 139210

$ws.Range("F450").Value = 87659

$ws.Range("F451").Value = 82384

$ws.Range("F452").Value = 72253

$ws.Range("F466").Value = 49170

$ws.Range("F467").Value = 49782

$ws.Range("F468").Value = 40189
$ws.Range("G468").Value = 44

$ws.Range("F469").Value = 38363

$ws.Range("F470").Value = 40422
$ws.Range("G470").Value = 41

# --- New rows appended at the bottom (dates 2021-06-18, -19, -20) ---
$newRows = @(
    @(471, 44365, 391297, 5671, 49, 12486, 54457, 37),
    @(472, 44366, 391325, 2537, 28, 12492, 40096, 18),
    @(473, 44367, 391326,  896,  1, 12496, 30101, 28)
)

foreach ($r in $newRows) {
    $row = $r[0]
    $ws.Cells.Item($row, 1).Value = $r[1]
    $ws.Cells.Item($row, 2).Value = $r[2]
    $ws.Cells.Item($row, 3).Value = $r[3]
    $ws.Cells.Item($row, 4).Value = $r[4]
    $ws.Cells.Item($row, 5).Value = $r[5]
    $ws.Cells.Item($row, 6).Value = $r[6]
    $ws.Cells.Item($row, 7).Value = $r[7]
}
